$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price (D) and Volume (E) cells to remain text so numeric-looking values
# (e.g. "1.034", "21.36") are not auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '27.510.46'
$ws.Range("E2").Value = '  +2.10%  '

$ws.Range("D3").Value = '1.846.96'
$ws.Range("E3").Value = '  +1.63%  '

$ws.Range("D4").Value = '1.034'
$ws.Range("E4").Value = '  +2.67%  '

$ws.Range("D5").Value = '320.87'
$ws.Range("E5").Value = '  +2.92%  '

$ws.Range("D6").Value = '1.029'
$ws.Range("E6").Value = '  +2.27%  '

$ws.Range("D7").Value = '0.4369'
$ws.Range("E7").Value = '  +1.85%  '

$ws.Range("D8").Value = '0.3762'
$ws.Range("E8").Value = '  +1.83%  '

$ws.Range("D9").Value = '0.07377'
$ws.Range("E9").Value = '  +1.85%  '

$ws.Range("D10").Value = '0.8721'
$ws.Range("E10").Value = '  +0.92%  '

$ws.Range("D11").Value = '21.36'
$ws.Range("E11").Value = '  +1.73%  '

$ws.Range("D12").Value = '1.862.28'
$ws.Range("E12").Value = '  -8.09%  '

$ws.Range("D13").Value = '5.497'
$ws.Range("E13").Value = '  +2.07%  '

$ws.Range("D14").Value = '6.664'
$ws.Range("E14").Value = '  +0.28%  '

$ws.Range("D15").Value = '0.07178'
$ws.Range("E15").Value = '  +3.63%  '

$ws.Range("D16").Value = '82.40'
$ws.Range("E16").Value = '  +2.03%  '

$ws.Range("E17").Value = '  +2.67%  '

$ws.Range("D18").Value = '0.000009032'
$ws.Range("E18").Value = '  +1.23%  '

$ws.Range("D19").Value = '1.028'
$ws.Range("E19").Value = '  +2.27%  '

$ws.Range("D20").Value = '15.37'
$ws.Range("E20").Value = '  +1.08%  '

$ws.Range("D21").Value = '27.533.05'
$ws.Range("E21").Value = '  +1.98%  '

$ws.Range("D22").Value = '5.232'
$ws.Range("E22").Value = '  +0.86%  '

$ws.Range("D23").Value = '11.31'
$ws.Range("E23").Value = '  +1.91%  '

$ws.Range("D24").Value = '2.071.17'
$ws.Range("E24").Value = '  -8.31%  '

$ws.Range("D25").Value = '157.39'
$ws.Range("E25").Value = '  +2.18%  '

$ws.Range("E26").Value = '  +1.76%  '

$ws.Range("E27").Value = '  +1.89%  '

$ws.Range("D28").Value = '5.252'
$ws.Range("E28").Value = '  +0.80%  '

$ws.Range("D29").Value = '1.946'
$ws.Range("E29").Value = '  +2.95%  '

$ws.Range("D30").Value = '116.59'
$ws.Range("E30").Value = '  +1.30%  '

$ws.Range("D31").Value = '0.09020'
$ws.Range("E31").Value = '  +0.61%  '

$ws.Range("D32").Value = '1.191'
$ws.Range("E32").Value = '  +2.22%  '

$ws.Range("D33").Value = '0.7588'
$ws.Range("E33").Value = '  +2.07%  '

$ws.Range("D34").Value = '4.490'
$ws.Range("E34").Value = '  +1.76%  '

$ws.Range("D35").Value = '2.874'
$ws.Range("E35").Value = '  +2.46%  '

$ws.Range("D36").Value = '1.029'
$ws.Range("E36").Value = '  +1.86%  '

$ws.Range("D37").Value = '1.147'
$ws.Range("E37").Value = '  +1.76%  '

$ws.Range("D38").Value = '0.01967'
$ws.Range("E38").Value = '  +2.07%  '

$ws.Range("D39").Value = '0.05274'
$ws.Range("E39").Value = '  +1.15%  '

$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = '0.5136'
$ws.Range("E40").Value = '  +1.24%  '

$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").Value = '2.797'
$ws.Range("E41").Value = '  +2.10%  '

$ws.Range("D42").Value = '0.1668'
$ws.Range("E42").Value = '  +1.45%  '

$ws.Range("D43").Value = '6.688'
$ws.Range("E43").Value = '  +4.08%  '

$ws.Range("D44").Value = '8.445'
$ws.Range("E44").Value = '  +2.75%  '

$ws.Range("D45").Value = '108.80'
$ws.Range("E45").Value = '  +1.64%  '

$ws.Range("D46").Value = '10.54'
$ws.Range("E46").Value = '  +1.53%  '

$ws.Range("D47").Value = '1.701'
$ws.Range("E47").Value = '  +2.65%  '

$ws.Range("D48").Value = '0.06398'
$ws.Range("E48").Value = '  +1.41%  '

$ws.Range("D49").Value = '0.4624'

$ws.Range("D50").Value = '1.852'
$ws.Range("E50").Value = '  +2.57%  '

$ws.Range("D51").Value = '39.09'
$ws.Range("E51").Value = '  +4.02%  '

# Restore default (General) formatting so the cell style matches the original
# workbook (no explicit style index), now that text values are safely stored.
$ws.Range("B2:E51").ClearFormats()
